# Bug-fix regression test fixture update:
# 1) Remove the stray _GoBack bookmark that Word leaves behind after edits.
# 2) Append a table containing an unnamed nested table, plus a trailing
#    paragraph, reproducing the "nested unnamed table replacement" test case.

$d = $word.ActiveDocument

# --- 1) Drop the _GoBack bookmark (hidden, but addressable by name). ---
try {
    $bm = $d.Bookmarks.Item("_GoBack")
    $bm.Delete()
} catch {
}

# --- 2) Append the new table (+ nested table) and a trailing paragraph. ---
# Word's Range model here always inserts new block content *before* the
# document's existing final (empty) paragraph when targeting Content.End,
# so we bundle a throw-away leading paragraph + the real table + the real
# trailing paragraph in one InsertXML call, then fold the now-redundant
# original trailing paragraph mark away again.
$beforeEnd = $d.Content.End

$endRange = $d.Range($beforeEnd, $beforeEnd)
$newBlockXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:tbl xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:tblPr><w:tblStyle w:val="TableGrid"/><w:tblW w:w="0" w:type="auto"/><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="4508"/><w:gridCol w:w="4508"/></w:tblGrid><w:tr><w:tc><w:tcPr><w:tcW w:w="4508" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>I DON</w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>''</w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>T WANT A T</w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>ABLE INSIDE A TABLE! :D</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4508" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="4508" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:instrText xml:space="preserve"> MERGEFIELD ANY_VALUE </w:instrText></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4508" w:type="dxa"/></w:tcPr><w:tbl><w:tblPr><w:tblStyle w:val="TableGrid"/><w:tblW w:w="0" w:type="auto"/><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="2141"/><w:gridCol w:w="2141"/></w:tblGrid><w:tr><w:tc><w:tcPr><w:tcW w:w="2141" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>Meehhh</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>, try to get rid of me</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2141" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:instrText xml:space="preserve"> MERGEFIELD ANY_VALUE </w:instrText></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="2141" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>well…</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2141" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">Done! </w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:instrText xml:space="preserve"> MERGEFIELD FORMAT_HINT_TABLE_REMOVE </w:instrText></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r></w:p></w:tc></w:tr></w:tbl><w:p><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr></w:p></w:tc></w:tr></w:tbl><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr></w:p>'
$endRange.InsertXML($newBlockXml)

$afterEnd = $d.Content.End
$mergeRange = $d.Range($afterEnd - 2, $afterEnd)
$mergeRange.Delete()

# --- 3) Make sure the final paragraph carries the en-GB language mark. ---
$finalRange = $d.Range($d.Content.End, $d.Content.End)
$finalRange.LanguageID = "en-GB"
